$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "0726563932"
$ws.Range("A2").Value = "Test Engineer Application"
$ws.Range("B2").Value = "Farouk"

$ws.Range("B2").Select()
